$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data (Step, X, Y) for rows 4 through 47, reflecting the updated
# mid-point / Bresenham iteration values from the commit.
$data = @(
    @(4,  2,  2,  1),
    @(5,  3,  3,  2),
    @(6,  4,  4,  3),
    @(7,  5,  5,  3),
    @(8,  6,  6,  4),
    @(9,  7,  7,  5),
    @(10, 8,  8,  5),
    @(11, 9,  9,  6),
    @(12, 10, 10, 6),
    @(13, 11, 11, 7),
    @(14, 12, 12, 8),
    @(15, 13, 13, 8),
    @(16, 14, 14, 9),
    @(17, 15, 15, 10),
    @(18, 16, 16, 10),
    @(19, 17, 17, 11),
    @(20, 18, 18, 12),
    @(21, 19, 19, 12),
    @(22, 20, 20, 13),
    @(23, 21, 21, 14),
    @(24, 22, 22, 14),
    @(25, 23, 23, 15),
    @(26, 24, 24, 15),
    @(27, 25, 25, 16),
    @(28, 26, 26, 17),
    @(29, 27, 27, 17),
    @(30, 28, 28, 18),
    @(31, 29, 29, 19),
    @(32, 30, 30, 19),
    @(33, 31, 31, 20),
    @(34, 32, 32, 21),
    @(35, 33, 33, 21),
    @(36, 34, 34, 22),
    @(37, 35, 35, 23),
    @(38, 36, 36, 23),
    @(39, 37, 37, 24),
    @(40, 38, 38, 24),
    @(41, 39, 39, 25),
    @(42, 40, 40, 26),
    @(43, 41, 41, 26),
    @(44, 42, 42, 27),
    @(45, 43, 43, 28),
    @(46, 44, 44, 28),
    @(47, 45, 45, 29)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
